$d = $word.ActiveDocument

$replacements = @(
    @{old="69×86="; new="67×84="},
    @{old="23×50="; new="43×73="},
    @{old="23×33="; new="71×21="},
    @{old="91×79="; new="55×75="},
    @{old="82×22="; new="64×74="},
    @{old="96×24="; new="81×76="},
    @{old="31×44="; new="53×97="},
    @{old="72×42="; new="21×13="},
    @{old="17×37="; new="54×56="},
    @{old="54×98="; new="35×36="},
    @{old="65×17="; new="23×15="},
    @{old="18×96="; new="61×22="},
    @{old="63×36="; new="74×11="},
    @{old="88×12="; new="49×28="},
    @{old="14×35="; new="66×55="},
    @{old="95×67="; new="56×14="},
    @{old="37×18="; new="42×95="},
    @{old="89×42="; new="47×82="},
    @{old="31×84="; new="90×23="},
    @{old="34×98="; new="37×29="},
    @{old="49×46="; new="36×13="},
    @{old="95×95="; new="43×71="},
    @{old="17×47="; new="77×48="},
    @{old="49×47="; new="89×99="},
    @{old="79×67="; new="97×17="}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
